$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.672.37"
$ws.Range("E2").Value = "  +0.51%  "
$ws.Range("D3").Value = "1.642.90"
$ws.Range("E3").Value = "  +0.86%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "215.04"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.02%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.504"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.49%  "
$ws.Range("E7").Value = "  -0.06%  "
$ws.Range("E8").Value = "  +0.52%  "
$ws.Range("E9").Value = "  +0.91%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.08"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.21%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0845"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.11%  "
$ws.Range("D12").Value = "1.870.32"
$ws.Range("E12").Value = "  +0.78%  "
$ws.Range("D13").Value = "1.643.10"
$ws.Range("E13").Value = "  +0.75%  "
$ws.Range("E14").Value = "  +1.76%  "
$ws.Range("E15").Value = "  +1.54%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "65.02"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.16%  "
$ws.Range("D17").Value = "26.684.40"
$ws.Range("E17").Value = "  +0.42%  "
$ws.Range("E18").Value = "  +0.57%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "215.59"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.40%  "
$ws.Range("E20").Value = "  -0.02%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.36"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.29%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.26"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.59%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.49"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.38%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.26"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +15.60%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "145.43"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.96%  "
$ws.Range("E26").Value = "  +0.07%  "
$ws.Range("E27").Value = "  -0.05%  "
$ws.Range("E28").Value = "  +4.21%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.70"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.27%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0517"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +2.26%  "
$ws.Range("E31").Value = "  +1.12%  "
$ws.Range("E32").Value = "  +2.33%  "
$ws.Range("E33").Value = "  +3.05%  "
$ws.Range("D34").Value = "1.281.70"
$ws.Range("E34").Value = "  +5.09%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.54"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +2.28%  "
$ws.Range("E36").Value = "  +1.27%  "
$ws.Range("E37").Value = "  +2.79%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.534"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +7.05%  "
$ws.Range("E39").Value = "  +4.41%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.01"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.02%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.817"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +2.96%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.26"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.12%  "
$ws.Range("E43").Value = "  +1.87%  "
$ws.Range("D44").Value = "1.781.39"
$ws.Range("E44").Value = "  +0.87%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "91.41"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.23%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "59.82"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +8.83%  "
$ws.Range("E47").Value = "  +1.91%  "
$ws.Range("E48").Value = "  +0.81%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "7.80"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.21%  "
$ws.Range("E50").Value = "  +2.24%  "
$ws.Range("E51").Value = "  -0.58%  "
